$d = $word.ActiveDocument
$paraCount = $d.Paragraphs.Count
Write-Output "ParaCount: $paraCount"

$r = $d.Paragraphs(3).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="573E31E2" w14:textId="69BF872F" w:rsidR="00066ACE" w:rsidRDefault="00066ACE"><w:r><w:t xml:space="preserve">Bijgevoegd een document met code voor het Globescope </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Assesment</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$r = $d.Paragraphs(4).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="242D2AFA" w14:textId="54D44F73" w:rsidR="00066ACE" w:rsidRDefault="00066ACE"><w:r><w:t xml:space="preserve">In het </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>assesment</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">wordt gevraagd om een programma te maken, waarmee minimaal 10 specifieke vragen  beantwoord kunnen worden. Dit programma is er op gericht om deze 10 vragen te beantwoorden, maar er kunnen ook </w:t></w:r><w:r><w:t xml:space="preserve">gelijkaardige </w:t></w:r><w:r><w:t>vragen mee beantwoord worden</w:t></w:r><w:r><w:t xml:space="preserve"> met iets andere input</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$r = $d.Paragraphs(6).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="680BED62" w14:textId="15D07B20" w:rsidR="009312D6" w:rsidRDefault="00066ACE" w:rsidP="00496E2B"><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>N</w:t></w:r><w:r><w:t xml:space="preserve">ode = plaats die </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>geidentificeerd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> wordt met een letter van het alfabet</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$r = $d.Paragraphs(7).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2C9EBCAB" w14:textId="23233263" w:rsidR="009312D6" w:rsidRDefault="00066ACE"><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>E</w:t></w:r><w:r><w:t>dge</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = verbinding tussen twee </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nodes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$r = $d.Paragraphs(8).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="60B8AD75" w14:textId="6F731482" w:rsidR="009312D6" w:rsidRDefault="00066ACE"><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>W</w:t></w:r><w:r><w:t>eight</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = afstand van de verbinding tussen twee </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nodes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$r = $d.Paragraphs(9).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="27E34A55" w14:textId="2820152E" w:rsidR="009312D6" w:rsidRDefault="00066ACE"><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>R</w:t></w:r><w:r><w:t xml:space="preserve">oute = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>georderde</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> lijst van afwisselend </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nodes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>edges</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, die de verbinding tussen een start node en een eind node vormen. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$r = $d.Paragraphs(10).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="55A756C6" w14:textId="382E84D4" w:rsidR="009312D6" w:rsidRDefault="00066ACE"><w:r><w:t xml:space="preserve">Dit vormde een uitgangspunt bij de opzet van de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>progamma</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> code.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$r = $d.Paragraphs(13).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4120CCDF" w14:textId="71BBB206" w:rsidR="009312D6" w:rsidRDefault="007E4B42"><w:r><w:t xml:space="preserve">De </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>progamma</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> code heb ik opgesteld in Python, omdat ik veel ervaring heb met deze taal. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$r = $d.Paragraphs(19).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="3DCE4DBA" w14:textId="6677C45D" w:rsidR="009312D6" w:rsidRDefault="007E4B42"><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>findroute</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">:  </w:t></w:r><w:r><w:t xml:space="preserve">voor het beantwoorden van </w:t></w:r><w:r><w:t>vragen voor rechtstreekse routes  tussen een begin- en een eindpunt.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$r = $d.Paragraphs(20).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2E88B329" w14:textId="63192C63" w:rsidR="009312D6" w:rsidRDefault="007E4B42"><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>f</w:t></w:r><w:r><w:t>indroute</w:t></w:r><w:r><w:t>_extend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t xml:space="preserve">voor het </w:t></w:r><w:r><w:t>beantwoord</w:t></w:r><w:r><w:t>en</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>van routes</w:t></w:r><w:r><w:t xml:space="preserve">, waarbij de route </w:t></w:r><w:r><w:t xml:space="preserve">kan </w:t></w:r><w:r><w:t>doorlo</w:t></w:r><w:r><w:t>pen na het bereiken van het eindpunt</w:t></w:r><w:r><w:t xml:space="preserve">.  </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$r = $d.Paragraphs(21).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4AEC330E" w14:textId="6C77F869" w:rsidR="00125924" w:rsidRDefault="00125924"><w:r><w:t xml:space="preserve">Door het aanroepen van deze twee functies </w:t></w:r><w:r><w:t xml:space="preserve">met de juiste input </w:t></w:r><w:r><w:t xml:space="preserve">kunnen de 10 vragen van het </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>assesment</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> beantwoord worden. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$r = $d.Paragraphs(23).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="399E9F1C" w14:textId="0C7A826B" w:rsidR="009312D6" w:rsidRDefault="00A006EC"><w:r><w:t xml:space="preserve">De hulpfuncties kunnen ook los van de functies </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>findroute</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>findroute_extend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> gebruikt worden, maar het toepassen ervan is niet uitgewerkt.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Naast deze functies zijn er ook classes.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$r = $d.Paragraphs(28).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4E51F840" w14:textId="21EF8D86" w:rsidR="009312D6" w:rsidRDefault="007E4B42"><w:r><w:t xml:space="preserve">Ik ben begonnen met classes voor een node- en een </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>edge</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-objecten te maken, “Node” en “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Edge</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">” class object. De </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>edges</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nodes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> worden samengesteld uit de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>graph</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. informatie, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>graph</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> informatie wordt omgezet naar node en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>edge</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-objecten via een aparte functie.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$r = $d.Paragraphs(29).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="39AD8A97" w14:textId="5A2D76EC" w:rsidR="009312D6" w:rsidRDefault="00496E2B"><w:r><w:t xml:space="preserve">In een later stadium </w:t></w:r><w:r><w:t xml:space="preserve">is een aparte klasse </w:t></w:r><w:r><w:t xml:space="preserve">gemaakt </w:t></w:r><w:r><w:t xml:space="preserve">voor een route, “Route”. Een route die de eind node bereikt had krijgt een True-waarde mee, routes die het eindpunt niet bereiken krijgen een waarde </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>False</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> mee</w:t></w:r><w:r><w:t xml:space="preserve"> voor het attribuut end.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$r = $d.Paragraphs(33).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="77BA1E94" w14:textId="77777777" w:rsidR="009312D6" w:rsidRDefault="007E4B42"><w:r><w:t xml:space="preserve">Met de test worden de uitkomsten van de 10 vragen van de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>assesment</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> getest.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$r = $d.Paragraphs(36).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="6915152E" w14:textId="79F180D1" w:rsidR="00B41A43" w:rsidRDefault="00A006EC" w:rsidP="00A006EC"><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>In functie “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>findroute</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">”: alleen </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>shortest</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> kan opgegeven worden voor </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>args</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$r = $d.Paragraphs(37).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2DE7649E" w14:textId="0C55CA79" w:rsidR="008917ED" w:rsidRDefault="00B25F08" w:rsidP="008917ED"><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Voor correcte werking van functie</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>findroute_extend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> moeten </w:t></w:r><w:r><w:t>Nodes</w:t></w:r><w:r><w:t xml:space="preserve"> uitwaartse verbinding met een node hebben.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$r = $d.Paragraphs(38).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="524534E2" w14:textId="656FC448" w:rsidR="008917ED" w:rsidRDefault="008917ED" w:rsidP="008917ED"><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Antwoorden op andere vragen dan de 10 uit het </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>assesment</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> zijn niet geverifieerd.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$r = $d.Paragraphs(39).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0BC9EBB4" w14:textId="3F121616" w:rsidR="00C57154" w:rsidRDefault="00C57154" w:rsidP="008917ED"><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Find_node</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>method</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> bij </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Edge</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> class: niet uitgewerkt voor </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>inout</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> =’out’</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$r = $d.Paragraphs(40).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="33E631C5" w14:textId="48D8BD79" w:rsidR="00C57154" w:rsidRDefault="00C57154" w:rsidP="008917ED"><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Find_node</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>method</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> bij </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Edge</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> class: mogelijkheid bestaat dat </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>noderes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> niet gedefinieerd wordt.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$r = $d.Paragraphs(41).Range
$r.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="3DE9B342" w14:textId="40CF16C5" w:rsidR="00E2122A" w:rsidRDefault="00E2122A" w:rsidP="008917ED"><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Findroute_extend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: no </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>such</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> route, is niet </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>geimplementeerd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
